$d = $word.ActiveDocument

# Locate the target paragraph: the long Greek description paragraph that
# currently references "Περσεύς" (Perseus) and should instead reference
# the "Αστερισμός Ωρίωνα" (Orion constellation).
$target = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*Συμμετέχετε*" -and $t -like "*παγκόσμια*" -and $t -like "*Περσεύς*") {
        $target = $p
        break
    }
}

$rng = $target.Range
# Exclude the trailing paragraph mark from the range we rewrite.
$rng.MoveEnd(1, -1) | Out-Null

# Clear all existing runs/formatting in the paragraph first so that the
# freshly-inserted text does not inherit any run-level formatting
# (matches the target: a single <w:r> with no <w:rPr>).
$rng.Text = ""

$rng2 = $target.Range
$rng2.MoveEnd(1, -1) | Out-Null

$newText = "Συμμετέχετε σε μία παγκόσμια καμπάνια για να παρατηρήσετε και να καταγράψετε τη φωτεινότητα των πιο αμυδρά ορατών άστρων σαν μέσο για την μέτρηση της Φωτορρύπανσης σε μία δεδομένη περιοχή. Με τον εντοπισμό και την παρατήρηση του  Αστερισμός Ωρίωνα στον νυχτερινό ουρανό καθώς και με την σύγκριση των ανωτέρω με τα διαγράμματα για τα μεγέθη των άστρων,  άνθρωποι από όλον τον κόσμο θα μάθουν πώς τα φώτα στην κοινότητά τους συμβάλλουν στην Φωτορρύπανση. Με την κατάθεση των πορισμάτων τους στην ιστοσελίδα θα δημιουργηθεί ένα αρχείο σχετικά με το τι μπορεί να δει κανείς στον νυχτερινό ουρανό."

$rng2.Text = $newText
